$d = $word.ActiveDocument

$replacements = @(
    @("154×2=", "357×4="),
    @("942×3=", "562×7="),
    @("383×3=", "961×4="),
    @("223×4=", "981×2="),
    @("845×7=", "811×4="),
    @("344×2=", "983×7="),
    @("953×9=", "486×7="),
    @("234×3=", "658×6="),
    @("367×2=", "846×3="),
    @("546×3=", "949×7="),
    @("310×7=", "610×7="),
    @("153×2=", "376×2="),
    @("103×2=", "763×8="),
    @("143×5=", "757×8="),
    @("841×6=", "734×6="),
    @("438×3=", "105×5="),
    @("373×5=", "253×2="),
    @("791×6=", "231×6="),
    @("570×5=", "739×2="),
    @("368×3=", "124×5="),
    @("833×3=", "462×8="),
    @("723×7=", "538×5="),
    @("586×3=", "412×6="),
    @("395×3=", "737×3="),
    @("458×8=", "134×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
